$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, pushing the existing row 92 (and
# everything below it) down by one.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new record.
$ws.Cells.Item(92, 1).Value = 3
$ws.Cells.Item(92, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(92, 3).Value = "Coquimbo"
$ws.Cells.Item(92, 4).Value = 44897
$ws.Cells.Item(92, 5).Value = 5
$ws.Cells.Item(92, 6).Value = 100112052
$ws.Cells.Item(92, 7).Value = "Albahaca"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 93
$ws.Cells.Item(92, 11).Value = 6500
$ws.Cells.Item(92, 12).Value = 7000
$ws.Cells.Item(92, 13).Value = 6796
$ws.Cells.Item(92, 14).Value = "$/docena de matas"
$ws.Cells.Item(92, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(92, 16).Value = 1133
$ws.Cells.Item(92, 17).Value = 6
$ws.Cells.Item(92, 18).Value = "Hortaliza"
